# The sheet gains a new top row (A1:K1) holding the plain column indices
# 0-10, while the original header row ("Lg.", "Threading", ...) is pushed
# down to row 2 (losing its bold/border header style, and three of its
# cells - H2, J2, K2 - becoming blank). Every data row shifts down by one
# (old row 2 -> row 3, ... old row 26 -> row 27), keeping its own content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row above the current row 1; this pushes the
# existing header row (and everything below it) down by one row.
$ws.Rows.Item(1).Insert()

# The row that used to be row 1 is now row 2 - copy its formatting
# (bold font + thin border) into the new row 1 so the header style moves
# with the numbering row instead of staying on the text-label row.
$ws.Range("A2:K2").Copy()
$ws.Range("A1:K1").PasteSpecial(-4122)

# Fill the new row 1 with the plain numeric column indices.
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
$ws.Range("K1").Value = 10

# Row 2 (the former row 1) drops the bold/bordered header style, matching
# the plain look the other label/data rows already have.
$ws.Range("A2:K2").Style = "Normal"

# Row 2 also loses the part number / thread-size / material text that
# used to live in H1/J1/K1 - those three cells become blank.
$ws.Range("H2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()
